# Rename the workbook's sheets to match the new naming scheme introduced
# by this commit ("clread build and updatedreadme"):
#   FileData -> SignUpTest
#   EditData -> SearchItem
$wb = $excel.ActiveWorkbook

$wb.Sheets.Item("FileData").Name = "SignUpTest"
$wb.Sheets.Item("EditData").Name = "SearchItem"
